$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 7 updates
$ws.Range("H7").Value = 2.8
$ws.Range("I7").Value = 3
$ws.Range("P7").Value = 1.5
$ws.Range("Q7").Value = 2.25
$ws.Range("R7").Value = 1.75
$ws.Range("S7").Value = 1.87
$ws.Range("T7").Value = 7.3
$ws.Range("U7").Value = 12.5
$ws.Range("V7").Value = 9.5
$ws.Range("W7").Value = 29
$ws.Range("Y7").Value = 32
$ws.Range("Z7").Value = 7.4
$ws.Range("AA7").Value = 5.5
$ws.Range("AB7").Value = 13
$ws.Range("AC7").Value = 65
$ws.Range("AD7").Value = 600
$ws.Range("AE7").Value = 8.25
$ws.Range("AF7").Value = 15.5
$ws.Range("AG7").Value = 10.5
$ws.Range("AH7").Value = 40
$ws.Range("AI7").Value = 28
$ws.Range("AJ7").Value = 35

# Row 9 updates
$ws.Range("G9").Value = 5.6
$ws.Range("H9").Value = 3.85
$ws.Range("M9").Value = 3.6
$ws.Range("N9").Value = 1.62
$ws.Range("R9").Value = 1.65
$ws.Range("S9").Value = 1.98
$ws.Range("T9").Value = 16
$ws.Range("U9").Value = 37
$ws.Range("V9").Value = 17.5
$ws.Range("W9").Value = 120
$ws.Range("X9").Value = 60
$ws.Range("Y9").Value = 50
$ws.Range("Z9").Value = 12
$ws.Range("AA9").Value = 7.7
$ws.Range("AE9").Value = 8
$ws.Range("AF9").Value = 8.25
$ws.Range("AG9").Value = 7.8
$ws.Range("AH9").Value = 12
$ws.Range("AI9").Value = 11.25
$ws.Range("AJ9").Value = 21
